# Re-order the header columns on row 2 of the Staging.Indicator template.
# The underlying shared-strings table is rebuilt by Excel to match the
# new left-to-right order of the header row, so we only need to rewrite
# the header cell values in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "IndicatorID",
    "BusinessKey",
    "IndicatorTypeBusinessKey",
    "OutcomeBusinessKey",
    "OutputBusinessKey",
    "ProgrammeBusinessKey",
    "ProjectBusinessKey",
    "SectorBusinessKey",
    "SubOutputBusinessKey",
    "SubSectorBusinessKey",
    "Baseline",
    "BaselineDate",
    "BaselineString",
    "Code",
    "LongName",
    "Notes",
    "ReleaseDate",
    "ReportingDate",
    "ShortName",
    "Target",
    "TargetDate",
    "TargetString",
    "TextDescription",
    "UnitOfMeasure"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newHeaders[$i]
}
